$d = $word.ActiveDocument

# The cell text is currently split across three runs (" Food Truck idea was
# pretty " / "stupid" / " and you're forced to sell it at a loss.") because a
# grammar-check proofing mark wraps "stupid". Re-typing the full sentence as a
# single Find & Replace merges it back into one run and drops the stray
# proofErr markers, matching the target edit.
$apostrophe = [char]0x2019
$findText = " Food Truck idea was pretty stupid and you" + $apostrophe + "re forced to sell it at a loss."

$range = $d.Content
$range.Find.ClearFormatting()
$range.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $findText, 2)
